$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10586
$ws.Range("H40").Value = 2246
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 2076.6667
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2076.6667
$ws.Range("M40").Value = -2325
$ws.Range("N40").Value = -2426.6667
$ws.Range("H98").Value = 2460.5
$ws.Range("I98").Value = 2065.5
$ws.Range("K98").Value = 2065.5
$ws.Range("M98").Value = -567.5
$ws.Range("H103").Value = 799.9
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 799.875
$ws.Range("K103").Value = 2400
$ws.Range("L103").Value = 2399.625
$ws.Range("M103").Value = -1814
$ws.Range("N103").Value = -3571.625
$ws.Range("H121").Value = 1076.2142
$ws.Range("J121").Value = 1076.2142
$ws.Range("L121").Value = 3228.6426
$ws.Range("N121").Value = -6722.642599999999
$ws.Range("H122").Value = 2460.5
$ws.Range("I122").Value = 2065.5
$ws.Range("K122").Value = 6196.5
$ws.Range("M122").Value = -3746.5
$ws.Range("H137").Value = 27305.41
$ws.Range("I137").Value = 763.3125
$ws.Range("J137").Value = 45769.477
$ws.Range("K137").Value = 2289.9375
$ws.Range("L137").Value = 137308.431
$ws.Range("M137").Value = 260.0625
$ws.Range("N137").Value = -142408.431

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3047.673
$ws.Range("I32").Value = 2172.25
$ws.Range("K32").Value = 2172.25
$ws.Range("M32").Value = -1885.25
$ws.Range("H61").Value = 4335
$ws.Range("I61").Value = 2954.4
$ws.Range("J61").Value = 5102
$ws.Range("K61").Value = 2954.4
$ws.Range("L61").Value = 5102
$ws.Range("M61").Value = -2742.4
$ws.Range("N61").Value = -5526
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 1810.3572
$ws.Range("I74").Value = 954.5
$ws.Range("K74").Value = 954.5
$ws.Range("M74").Value = -80.5
$ws.Range("H77").Value = 1810.3572
$ws.Range("I77").Value = 954.5
$ws.Range("K77").Value = 4772.5
$ws.Range("M77").Value = -404.5
$ws.Range("H132").Value = 3303.7334
$ws.Range("I132").Value = 3231.818
$ws.Range("J132").Value = 3501.5
$ws.Range("K132").Value = 9695.454000000002
$ws.Range("L132").Value = 10504.5
$ws.Range("M132").Value = -7165.454000000002
$ws.Range("N132").Value = -15564.5
$ws.Range("H136").Value = 4335
$ws.Range("I136").Value = 2954.4
$ws.Range("J136").Value = 5102
$ws.Range("K136").Value = 8863.200000000001
$ws.Range("L136").Value = 15306
$ws.Range("M136").Value = -6313.200000000001
$ws.Range("N136").Value = -20406

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1555
$ws.Range("I5").Value = 1555
$ws.Range("K5").Value = 1555
$ws.Range("M5").Value = -1442
$ws.Range("H134").Value = 3034
$ws.Range("I134").Value = 2896
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 8688
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -6153
$ws.Range("N134").Value = -17070

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1265.5435
$ws.Range("I31").Value = 984.3570999999999
$ws.Range("J31").Value = 1388.5625
$ws.Range("K31").Value = 984.3570999999999
$ws.Range("L31").Value = 1388.5625
$ws.Range("M31").Value = -689.3570999999999
$ws.Range("N31").Value = -1978.5625
$ws.Range("H34").Value = 1265.5435
$ws.Range("I34").Value = 984.3570999999999
$ws.Range("J34").Value = 1388.5625
$ws.Range("K34").Value = 984.3570999999999
$ws.Range("L34").Value = 1388.5625
$ws.Range("M34").Value = -782.3570999999999
$ws.Range("N34").Value = -1792.5625
$ws.Range("H99").Value = 1002247.7
$ws.Range("I99").Value = 1430282.4
$ws.Range("K99").Value = 1430282.4
$ws.Range("M99").Value = -1428784.4
$ws.Range("H126").Value = 1002247.7
$ws.Range("I126").Value = 1430282.4
$ws.Range("K126").Value = 4290847.199999999
$ws.Range("M126").Value = -4288377.199999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1930.4706
$ws.Range("I68").Value = 837.6
$ws.Range("K68").Value = 2512.8
$ws.Range("M68").Value = -1701.8
$ws.Range("H71").Value = 1930.4706
$ws.Range("I71").Value = 837.6
$ws.Range("K71").Value = 7538.400000000001
$ws.Range("M71").Value = -3482.400000000001
$ws.Range("H86").Value = 2021.8
$ws.Range("I86").Value = 2152.25
$ws.Range("K86").Value = 6456.75
$ws.Range("M86").Value = -5270.75
$ws.Range("H89").Value = 2021.8
$ws.Range("I89").Value = 2152.25
$ws.Range("K89").Value = 19370.25
$ws.Range("M89").Value = -13442.25
$ws.Range("H98").Value = 466.9
$ws.Range("I98").Value = 466
$ws.Range("K98").Value = 1398
$ws.Range("M98").Value = 100
$ws.Range("H131").Value = 33382944
$ws.Range("J131").Value = 67481.63
$ws.Range("L131").Value = 202444.89
$ws.Range("N131").Value = -212524.89
$ws.Range("H140").Value = 3018.543
$ws.Range("I140").Value = 576.3333
$ws.Range("J140").Value = 5604.4116
$ws.Range("K140").Value = 1728.9999
$ws.Range("L140").Value = 16813.2348
$ws.Range("M140").Value = 3451.0001
$ws.Range("N140").Value = -27173.2348

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2349.1428
$ws.Range("I102").Value = 2376.4614
$ws.Range("K102").Value = 2376.4614
$ws.Range("M102").Value = -754.4614000000001
$ws.Range("H122").Value = 3297.8
$ws.Range("I122").Value = 1504.75
$ws.Range("J122").Value = 4493.1665
$ws.Range("K122").Value = 4514.25
$ws.Range("L122").Value = 13479.4995
$ws.Range("M122").Value = -2064.25
$ws.Range("N122").Value = -18379.4995

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7678.55
$ws.Range("I122").Value = 5684.0713
$ws.Range("J122").Value = 12332.333
$ws.Range("K122").Value = 17052.2139
$ws.Range("L122").Value = 36996.999
$ws.Range("M122").Value = -14602.2139
$ws.Range("N122").Value = -41896.999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4499.2856
$ws.Range("I132").Value = 3400
$ws.Range("J132").Value = 4939
$ws.Range("K132").Value = 10200
$ws.Range("L132").Value = 14817
$ws.Range("M132").Value = -7670
$ws.Range("N132").Value = -19877
$ws.Range("H136").Value = 19843960
$ws.Range("I136").Value = 32682146
$ws.Range("K136").Value = 98046438
$ws.Range("M136").Value = -98043888
